$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.094.09'
$ws.Range('E2').Value = '  +2.73%  '

$ws.Range('D3').Value = '2.966.79'
$ws.Range('E3').Value = '  +1.02%  '

$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('D5').Value = '595.67'
$ws.Range('E5').Value = '  -0.06%  '

$ws.Range('D6').Value = '147.72'
$ws.Range('E6').Value = '  +1.71%  '

$ws.Range('E7').Value = '  -0.08%  '

$ws.Range('D8').Value = '2.964.04'
$ws.Range('E8').Value = '  +1.03%  '

$ws.Range('D9').Value = '0.509'
$ws.Range('E9').Value = '  +0.80%  '

$ws.Range('D10').Value = '7.26'
$ws.Range('E10').Value = '  +3.76%  '

$ws.Range('D11').Value = '0.153'
$ws.Range('E11').Value = '  +6.60%  '

$ws.Range('E12').Value = '  +1.02%  '

$ws.Range('E13').Value = '  +6.19%  '

$ws.Range('D14').Value = '33.35'
$ws.Range('E14').Value = '  -1.04%  '

$ws.Range('E15').Value = '  -0.51%  '

$ws.Range('D16').Value = '3.461.07'
$ws.Range('E16').Value = '  +0.22%  '

$ws.Range('D17').Value = '62.940.77'
$ws.Range('E17').Value = '  +2.63%  '

$ws.Range('D18').Value = '6.77'
$ws.Range('E18').Value = '  +0.32%  '

$ws.Range('D19').Value = '2.967.58'
$ws.Range('E19').Value = '  +1.14%  '

$ws.Range('D20').Value = '445.15'
$ws.Range('E20').Value = '  +2.54%  '

$ws.Range('D21').Value = '13.53'
$ws.Range('E21').Value = '  +0.00%  '

$ws.Range('D22').Value = '0.671'
$ws.Range('E22').Value = '  -1.49%  '

$ws.Range('D23').Value = '7.13'
$ws.Range('E23').Value = '  -0.78%  '

$ws.Range('D24').Value = '11.30'
$ws.Range('E24').Value = '  +3.10%  '

$ws.Range('D25').Value = '81.77'
$ws.Range('E25').Value = '  -0.27%  '

$ws.Range('D26').Value = '2.16'
$ws.Range('E26').Value = '  -2.95%  '

$ws.Range('D27').Value = '11.94'
$ws.Range('E27').Value = '  +0.46%  '

$ws.Range('E28').Value = '  +0.01%  '

$ws.Range('D29').Value = '7.35'
$ws.Range('E29').Value = '  +4.80%  '

$ws.Range('D30').Value = '2.64'
$ws.Range('E30').Value = '  +0.61%  '

$ws.Range('E31').Value = '  -3.45%  '

$ws.Range('D32').Value = '0.0₃0977'
$ws.Range('E32').Value = '  +10.08%  '

$ws.Range('D33').Value = '26.64'
$ws.Range('E33').Value = '  -0.64%  '

$ws.Range('E34').Value = '  -1.79%  '

$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.09%  '

$ws.Range('D36').Value = '0.997'
$ws.Range('E36').Value = '  -1.72%  '

$ws.Range('E37').Value = '  +3.82%  '

$ws.Range('E38').Value = '  -0.06%  '

$ws.Range('D39').Value = '2.08'
$ws.Range('E39').Value = '  +2.17%  '

$ws.Range('D40').Value = '49.58'
$ws.Range('E40').Value = '  -0.58%  '

$ws.Range('D41').Value = '8.56'
$ws.Range('E41').Value = '  -0.99%  '

$ws.Range('E42').Value = '  -4.37%  '

$ws.Range('D43').Value = '0.285'
$ws.Range('E43').Value = '  +0.08%  '

$ws.Range('D44').Value = '40.94'
$ws.Range('E44').Value = '  -3.98%  '

$ws.Range('D45').Value = '2.718.71'
$ws.Range('E45').Value = '  +0.08%  '

$ws.Range('D46').Value = '134.43'
$ws.Range('E46').Value = '  +0.69%  '

$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').Value = '0.0340'
$ws.Range('E47').Value = '  -2.83%  '

$ws.Range('B48').Value = 'Bittensor'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D48').Value = '364.78'
$ws.Range('E48').Value = '  -2.21%  '

$ws.Range('E50').Value = '  -0.52%  '

$ws.Range('D51').Value = '23.02'
$ws.Range('E51').Value = '  -4.14%  '
